$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 209, shifting existing rows 209.. down to 210..
$ws.Rows.Item(209).Insert()

# Populate the newly inserted row 209 with the new record's data.
$ws.Range("A209").Value = 4
$ws.Range("B209").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C209").Value = "Los Lagos"
$ws.Range("D209").Value = 44876
$ws.Range("E209").Value = 10
$ws.Range("F209").Value = "Fruta"
$ws.Range("G209").Value = 100108
$ws.Range("H209").Value = "Tropicales y subtropicales"
$ws.Range("I209").Value = 100108002
$ws.Range("J209").Value = "Mango"
$ws.Range("K209").Value = "Sin especificar"
$ws.Range("L209").Value = "Primera"
$ws.Range("M209").Value = 200
$ws.Range("N209").Value = 8000
$ws.Range("O209").Value = 9000
$ws.Range("P209").Value = 8500
$ws.Range("Q209").Value = "$/bandeja 4 kilos"
$ws.Range("R209").Value = "Brasil"
$ws.Range("S209").Value = 2125
$ws.Range("T209").Value = 4
